# "Corrections has been made"
#
# 1) RN001 label: make the "RN001:" label bold, keep the following space
#    and the rest of the sentence in regular weight. This turns the single
#    run "RN001: " into a bold run "RN001:" followed by a plain run " ".
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("RN001:", $true, $true, $false, $false, $false, `
                            $true, 1, $false, $null, 0)
if ($found) {
    $rng.Bold = 1
}

# 2) Normalize the built-in heading style display names so they read
#    "Heading 1".."Heading 6" instead of the lower-cased "heading 1".."heading 6".
$headingMap = @{
    "Ttulo1" = "Heading 1"
    "Ttulo2" = "Heading 2"
    "Ttulo3" = "Heading 3"
    "Ttulo4" = "Heading 4"
    "Ttulo5" = "Heading 5"
    "Ttulo6" = "Heading 6"
}
foreach ($styleId in $headingMap.Keys) {
    try {
        $style = $d.Styles($styleId)
        if ($style -ne $null) {
            $style.NameLocal = $headingMap[$styleId]
        }
    } catch {
    }
}
